# Update countries & provincias Spain
# Refresh the COVID-19 country stats (paises.xlsx / sheet "Pais"):
#  - bump the "Datos actualizados" timestamp in A1
#  - update the per-country numeric columns (B..H) that changed between
#    the 00:42 and 01:59 data pulls
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp in header cell A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 01:59"

# Row 4
$ws.Range("B4").Value = 2681312
$ws.Range("C4").Value = 44235
$ws.Range("D4").Value = 1104679
$ws.Range("E4").Value = 1447861
$ws.Range("G4").Value = 335
$ws.Range("H4").Value = 128772

# Row 5
$ws.Range("B5").Value = 1370488
$ws.Range("C5").Value = 25234
$ws.Range("E5").Value = 554641
$ws.Range("G5").Value = 727
$ws.Range("H5").Value = 58385

# Row 22
$ws.Range("B22").Value = 103918
$ws.Range("C22").Value = 668
$ws.Range("D22").Value = 67178
$ws.Range("E22").Value = 28174

# Row 28
$ws.Range("B28").Value = 62268
$ws.Range("C28").Value = 2335
$ws.Range("D28").Value = 21138
$ws.Range("E28").Value = 39850
$ws.Range("G28").Value = 48
$ws.Range("H28").Value = 1280

# Row 29
$ws.Range("B29").Value = 61790
$ws.Range("C29").Value = 315
$ws.Range("D29").Value = 45213
$ws.Range("E29").Value = 16190
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 387

# Row 30
$ws.Range("B30").Value = 61361
$ws.Range("C30").Value = 66
$ws.Range("D30").Value = 16941
$ws.Range("E30").Value = 34688
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 9732

# Row 85
$ws.Range("B85").Value = 5394
$ws.Range("C85").Value = 185
$ws.Range("D85").Value = 2420
$ws.Range("E85").Value = 2932
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 42

# Row 86
$ws.Range("B86").Value = 5351
$ws.Range("C86").Value = 9
$ws.Range("D86").Value = 4296
$ws.Range("E86").Value = 1024
$ws.Range("H86").Value = 31

# Row 87
$ws.Range("B87").Value = 5297
$ws.Range("D87").Value = 1649
$ws.Range("E87").Value = 3604
$ws.Range("H87").Value = 44

# Row 96
$ws.Range("B96").Value = 3531
$ws.Range("D96").Value = 754
$ws.Range("E96").Value = 2732

# Row 100
$ws.Range("B100").Value = 2904
$ws.Range("C100").Value = 10
$ws.Range("D100").Value = 910
$ws.Range("E100").Value = 1904

# Row 115
$ws.Range("B115").Value = 1840
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 1818

# Row 133
$ws.Range("B133").Value = 1075
$ws.Range("C133").Value = 1
$ws.Range("E133").Value = 69

# Row 137
$ws.Range("B137").Value = 932
$ws.Range("C137").Value = 3
$ws.Range("D137").Value = 822
$ws.Range("E137").Value = 83

# Row 143
$ws.Range("B143").Value = 802
$ws.Range("C143").Value = 40
$ws.Range("D143").Value = 206
$ws.Range("E143").Value = 573
$ws.Range("G143").Value = 2
$ws.Range("H143").Value = 23

# Row 144
$ws.Range("B144").Value = 795
$ws.Range("C144").Value = 14
$ws.Range("D144").Value = 380
$ws.Range("E144").Value = 404
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 11

# Row 145
$ws.Range("B145").Value = 770
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 312
$ws.Range("E145").Value = 422
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 36

# Row 146
$ws.Range("D146").Value = 235
$ws.Range("E146").Value = 465

# Row 151
$ws.Range("B151").Value = 643
$ws.Range("C151").Value = 1
$ws.Range("E151").Value = 228

# Row 152
$ws.Range("B152").Value = 574
$ws.Range("C152").Value = 7
$ws.Range("D152").Value = 152
$ws.Range("E152").Value = 415
$ws.Range("G152").Value = 1
$ws.Range("H152").Value = 7

# Row 155
$ws.Range("C155").Value = 11
$ws.Range("D155").Value = 217
$ws.Range("E155").Value = 271
$ws.Range("G155").Value = 2
$ws.Range("H155").Value = 13

# Row 156
$ws.Range("B156").Value = 501
$ws.Range("C156").Value = 20
$ws.Range("D156").Value = 315
$ws.Range("E156").Value = 175
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 11

# Row 158
$ws.Range("D158").Value = 335
$ws.Range("E158").Value = 20
